$d = $word.ActiveDocument

# The document contains four inline "spacer" pictures (1x1 px placeholder
# images) that illustrate various height-control diagrams. Each is being
# replaced by a plain hyperlink run whose display text is the URL of the
# actual illustration image hosted on ura.gov.sg, and whose target address
# is that same URL.

$replacements = @(
    "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Flats-Condominiums/F05_Building_Height.jpg?h=100%25&w=100%25",
    "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Flats-Condominiums/Singapore-Botanic-Gardens.jpg?h=100%25&w=100%25",
    "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Flats-Condominiums/F06_Floor_to_Floor_Height.jpg?h=100%25&w=100%25",
    "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Commercial/C04_Additional_Height_for_Sky_Terrace_Floors.jpg?h=100%25&w=100%25"
)

# Walk the InlineShapes from last to first so that earlier shapes' Range
# positions are not invalidated by edits made to later ones.
$shapes = $d.InlineShapes
$count = $shapes.Count
for ($i = $count; $i -ge 1; $i--) {
    $shape = $shapes.Item($i)
    $url = $replacements[$i - 1]
    $r = $shape.Range
    $r.Text = $url
    $d.Hyperlinks.Add($r, $url) | Out-Null
}

Write-Output "Replaced $count inline images with hyperlinks"
